# Insert a new weekly price record as row 243, pushing the existing
# rows 243:333 down to 244:334 (dimension grows from A1:R333 to A1:R334).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 243:333 down by one row.
$ws.Rows.Item(243).Insert()

# Populate the newly inserted row 243 with the new record.
$ws.Cells.Item(243, 1).Value2  = 9
$ws.Cells.Item(243, 2).Value2  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(243, 3).Value2  = "Metropolitana"
$ws.Cells.Item(243, 4).Value2  = 44627
$ws.Cells.Item(243, 5).Value2  = 13
$ws.Cells.Item(243, 6).Value2  = 100112032
$ws.Cells.Item(243, 7).Value2  = "Zapallo italiano"
$ws.Cells.Item(243, 8).Value2  = "Sin especificar"
$ws.Cells.Item(243, 9).Value2  = "Primera"
$ws.Cells.Item(243, 10).Value2 = 79
$ws.Cells.Item(243, 11).Value2 = 14000
$ws.Cells.Item(243, 12).Value2 = 15000
$ws.Cells.Item(243, 13).Value2 = 14494
$ws.Cells.Item(243, 14).Value2 = "$/caja 50 unidades"
$ws.Cells.Item(243, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(243, 16).Value2 = 290
$ws.Cells.Item(243, 17).Value2 = 50
$ws.Cells.Item(243, 18).Value2 = "Hortaliza"
